# Weekly update: insert the new week's record for
# "Vega Modelo de Temuco - Bruselas (repollito)" at row 35,
# shifting the existing rows 35:50 down to 36:51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (shifts rows 35..50 down to 36..51,
# carrying formatting, including the date style on column D).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with this week's data.
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 44455
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = 100112035
$ws.Cells.Item(35, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 20
$ws.Cells.Item(35, 11).Value = 25000
$ws.Cells.Item(35, 12).Value = 25000
$ws.Cells.Item(35, 13).Value = 25000
$ws.Cells.Item(35, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(35, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(35, 16).Value = 2500
$ws.Cells.Item(35, 17).Value = 10
$ws.Cells.Item(35, 18).Value = "Hortaliza"
